$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999999549822216,
    0.99555708568604973,
    0.98017271900621883,
    0.97466217561852275,
    0.96970331376184427,
    0.94633035167620727,
    0.93973076519590659,
    0.93059689826335945,
    0.91884633512850022,
    0.90801439501239489,
    0.90641783656284325,
    0.90363730548318366,
    0.89234782041183802,
    0.88817990781193323,
    0.88558813337578335,
    0.88308144195698168,
    0.87937330729540542,
    0.87826437583905992,
    0.99450382753497402,
    0.98738671456911464,
    0.98598821503835055,
    0.98472370931551056,
    0.97313758065203704,
    0.96011636183913396,
    0.9536593052017569,
    0.9450813655742254,
    0.94117841539685865,
    0.92435850256244856,
    0.91252938678847195,
    0.90756462711661912,
    0.90434245386450884,
    0.89978920820123276,
    0.89553359975105895
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
